# Applies the commit's change: wraps eight trainee names in the
# "Trainees" column of the table with spell-check proofErr markers
# (w:proofErr type="spellStart"/"spellEnd"), and fills in the
# previously-empty Trainees cell for row 24 (ELK24) with "Arun".
#
# The runtime's Range.InsertXML() only preserves non-w:p child markup
# (like w:proofErr) when the target Range spans an entire paragraph
# (start through its end-of-paragraph mark); it is then replaced in
# place. So each edit below selects the whole paragraph via the
# table's Cell(row, col).Range.Paragraphs(1).Range and replaces it
# with a reconstructed <w:p> carrying the same w14:paraId/w14:textId/
# w:rsidR/w:rsidRDefault attributes (and rPr/pPr) as the original,
# plus the new markup.
#
# NOTE: this interpreter's functions only bind POSITIONAL arguments
# (named "-Param value" style binding does not populate $Param), so
# the helper below takes plain positional parameters.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$wNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"
$w14Ns = "xmlns:w14='http://schemas.microsoft.com/office/word/2010/wordml'"

function Add-SpellCheckMarks {
    param($Row, $Name, $ParaId, $TextId, $RsidRDefault)

    $cell = $t.Cell($Row, 4)
    $para = $cell.Range.Paragraphs.Item(1)
    $rng = $para.Range

    $xml = "<w:p $wNs $w14Ns w14:paraId='$ParaId' w14:textId='$TextId' w:rsidR='00BF39C6' w:rsidRDefault='$RsidRDefault'>" +
           "<w:pPr><w:rPr><w:color w:val='000000'/></w:rPr></w:pPr>" +
           "<w:proofErr w:type='spellStart'/>" +
           "<w:r><w:rPr><w:color w:val='000000'/></w:rPr><w:t>$Name</w:t></w:r>" +
           "<w:proofErr w:type='spellEnd'/>" +
           "</w:p>"

    $rng.InsertXML($xml)
}

Add-SpellCheckMarks 12 "Dilip"       "23658185" "2F7A06B7" "00BF39C6"
Add-SpellCheckMarks 14 "Sugunan"     "1EFB2520" "43F0653B" "00BF39C6"
Add-SpellCheckMarks 15 "Dananjay"    "648A074D" "3F5404A1" "00BF39C6"
Add-SpellCheckMarks 16 "Vajira"      "25804DFE" "4D27B7B5" "00BF39C6"
Add-SpellCheckMarks 18 "Nisansala"   "10D310A0" "37C66E25" "00BF39C6"
Add-SpellCheckMarks 19 "Keerthana"   "63167E52" "142ADC1C" "00BF39C6"
Add-SpellCheckMarks 23 "charumathi"  "1EF54297" "45445B7C" "00923868"
Add-SpellCheckMarks 25 "Sabarivasan" "46A6C52D" "6F42A4ED" "00923868"

# Row 26 (SI.NO 24 / ELK24) Trainees cell was empty; fill it with "Arun".
$arunCell = $t.Cell(26, 4)
$arunPara = $arunCell.Range.Paragraphs.Item(1)
$arunRng = $arunPara.Range

$arunXml = "<w:p $wNs $w14Ns w14:paraId='6F62CE3E' w14:textId='77777777' w:rsidR='00BF39C6' w:rsidRDefault='00BF39C6'>" +
           "<w:pPr><w:rPr><w:color w:val='000000'/></w:rPr></w:pPr>" +
           "<w:r><w:rPr><w:color w:val='000000'/></w:rPr><w:t>Arun</w:t></w:r>" +
           "</w:p>"

$arunRng.InsertXML($arunXml)

Write-Host "Applied proofErr wraps and Arun fill-in."
